$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:F5 as described in the diff
$data = @{
    2 = @(1.175720741260731, 1.175720741260731, 2.00046158807684, 1.414376748987638, 0.8158993420687196)
    3 = @(1.179154357706104, 1.179154357706104, 1.839172490612837, 1.356160938315522, 0.7061377669450614)
    4 = @(0.9948462645758421, 0.9948462645758421, 1.24218013291981, 1.114531351250296, 0.5504118924361679)
    5 = @(1.219879926150666, 1.219879926150666, 1.817306308435718, 1.348075038132417, 0.8114176165333894)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
}
